# feat: add 2022-Q3 data
#
# The workbook starts with two sheets: "总计" (summary) and "2022-Q1" (the
# only quarter of fund-holding detail recorded so far).
#
# The edit: a new quarter ("2022-Q3") of fund-holding detail is recorded.
# The sheet that used to hold the (only) quarter of detail is renamed to
# "2022-Q3" and its content is replaced by the new quarter's data; a fresh
# sheet named "2022-Q1" is inserted right after it, re-creating the
# original "2022-Q1" detail so that history is preserved. The "总计"
# summary sheet gains a row for the new quarter while keeping the old one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet ("总计"): insert the new quarter above the existing row
#    by pushing the old "2022-Q1" summary row down to row 3, then writing
#    the "2022-Q3" summary into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 3 is new - give A3 the same "index" look as A2 (bold/bordered style)
# before filling in the values that used to live in row 2.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.27

# Row 2 now becomes the newest quarter, "2022-Q3".
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.25

# ---------------------------------------------------------------------
# 2. Detail sheets: duplicate the existing "2022-Q1" sheet so the old data
#    keeps living on its own tab, then rename the original tab to
#    "2022-Q3" and overwrite its contents with the new quarter's holdings.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")

# Duplicate "2022-Q1" right after itself - the copy keeps its data/format
# and becomes the new "2022-Q1" tab once we rename the original below.
$q1.Copy($null, $q1)

$q3 = $q1
$q3.Name = "2022-Q3"
$newQ1 = $wb.Worksheets.Item("2022-Q1 (2)")
$newQ1.Name = "2022-Q1"

# Wipe the old "2022-Q1" holdings out of what is now the "2022-Q3" tab.
$q3.Cells.Clear()

# Re-create the header row, copying the bold/bordered header format that
# is already used on the "总计" sheet so no new style gets introduced.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$summary.Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)

# Fund codes must stay text - "001628" would lose its leading zero if
# Excel were left to treat it as a number. Same for the size/weight
# columns, which already read as text on the other quarter's sheet.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

$data = @(
    @(0, "513360", "博时中证全球中国教育主题ETF（QDII）", "4.81", "99.43", "2.56", "0.1231", 10),
    @(1, "001628", "招商体育文化休闲股票A",                 "2.23", "92.42", "5.14", "0.1146", 3),
    @(2, "015395", "招商体育文化休闲股票C",                 "0.25", "92.42", "5.14", "0.0128", 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Leave the re-created "2022-Q1" sheet as the active tab, same as before
# the edit.
$newQ1.Activate()
$newQ1.Range("A1").Select()
